$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.857.71"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.177.71"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("E4").Value = "  +0.05%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "571.43"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "171.99"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -3.35%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.600"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "3.175.81"
$ws.Range("E9").Value = "  -4.93%  "
$ws.Range("E10").Value = "  -3.36%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "6.59"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -4.38%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.393"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -5.05%  "
$ws.Range("D13").Value = "3.730.80"
$ws.Range("E13").Value = "  -4.83%  "
$ws.Range("E14").Value = "  +1.14%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "27.43"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -4.70%  "
$ws.Range("D16").Value = "65.818.26"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18").Value = "3.176.66"
$ws.Range("E18").Value = "  -5.08%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "5.72"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "12.91"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "360.12"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("E22").Value = "  -2.00%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "69.38"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -2.99%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.496"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -4.93%  "
$ws.Range("D26").Value = "3.306.24"
$ws.Range("E26").Value = "  -5.15%  "
$ws.Range("E27").Value = "  -6.43%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "9.86"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -0.12%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "1.93"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -1.96%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "5.42"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -4.21%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "22.03"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -3.76%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.21"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "6.63"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -3.57%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "159.74"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -3.05%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.837"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "1.81"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "26.44"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -3.35%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "2.51"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "2.657.65"
$ws.Range("E43").Value = "  -1.59%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "6.13"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("E45").Value = "  -2.65%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "39.67"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0660"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "330.20"
$r.Style = "Normal"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "24.19"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  -1.46%  "
